$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.794.86"
$ws.Range("E2").Value = "  +2.51%  "

$ws.Range("D3").Value = "2.609.68"
$ws.Range("E3").Value = "  +0.84%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'567.10"
$ws.Range("E5").Value = "  -0.85%  "

$ws.Range("D6").Value = "'142.67"
$ws.Range("E6").Value = "  -0.74%  "

$ws.Range("E7").Value = "  -0.26%  "

$ws.Range("D8").Value = "'0.603"
$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").Value = "2.631.47"
$ws.Range("E9").Value = "  +1.32%  "

$ws.Range("D10").Value = "'6.73"
$ws.Range("E10").Value = "  +1.01%  "

$ws.Range("E11").Value = "  +1.82%  "

$ws.Range("E12").Value = "  +1.47%  "

$ws.Range("D13").Value = "'0.375"
$ws.Range("E13").Value = "  +8.15%  "

$ws.Range("D14").Value = "3.073.33"
$ws.Range("E14").Value = "  +1.03%  "

$ws.Range("D15").Value = "60.711.09"
$ws.Range("E15").Value = "  +2.31%  "

$ws.Range("D16").Value = "'23.47"
$ws.Range("E16").Value = "  +4.13%  "

$ws.Range("D18").Value = "2.614.94"
$ws.Range("E18").Value = "  +0.94%  "

$ws.Range("D19").Value = "'4.68"
$ws.Range("E19").Value = "  +3.17%  "

$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'348.41"
$ws.Range("E20").Value = "  +2.88%  "

$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").Value = "'10.95"
$ws.Range("E21").Value = "  +6.93%  "

$ws.Range("D22").Value = "'7.02"
$ws.Range("E22").Value = "  +13.16%  "

$ws.Range("E23").Value = "  +0.11%  "

$ws.Range("D24").Value = "'0.519"
$ws.Range("E24").Value = "  +14.45%  "

$ws.Range("D25").Value = "'63.42"
$ws.Range("E25").Value = "  -1.64%  "

$ws.Range("D26").Value = "'0.995"
$ws.Range("E26").Value = "  -0.45%  "

$ws.Range("D27").Value = "'0.162"
$ws.Range("E27").Value = "  -0.28%  "

$ws.Range("E28").Value = "  +6.23%  "

$ws.Range("D29").Value = "0.0₃0794"
$ws.Range("E29").Value = "  +1.50%  "

$ws.Range("D30").Value = "'1.78"
$ws.Range("E30").Value = "  +5.26%  "

$ws.Range("E31").Value = "  -0.10%  "

$ws.Range("D32").Value = "'6.27"
$ws.Range("E32").Value = "  +3.13%  "

$ws.Range("D33").Value = "'161.42"
$ws.Range("E33").Value = "  +1.52%  "

$ws.Range("D34").Value = "'19.60"
$ws.Range("E34").Value = "  +2.74%  "

$ws.Range("E35").Value = "  +4.42%  "

$ws.Range("D36").Value = "'0.961"
$ws.Range("E36").Value = "  +8.63%  "

$ws.Range("E37").Value = "  +4.49%  "

$ws.Range("E38").Value = "  +5.70%  "

$ws.Range("D39").Value = "'37.68"
$ws.Range("E39").Value = "  +1.60%  "

$ws.Range("E40").Value = "  -2.85%  "

$ws.Range("E41").Value = "  +3.90%  "

$ws.Range("D42").Value = "'302.99"
$ws.Range("E42").Value = "  +2.70%  "

$ws.Range("D43").Value = "'141.35"
$ws.Range("E43").Value = "  +13.69%  "

$ws.Range("D44").Value = "'0.994"
$ws.Range("E44").Value = "  -0.40%  "

$ws.Range("D45").Value = "'0.0984"
$ws.Range("E45").Value = "  +0.64%  "

$ws.Range("D46").Value = "'0.604"
$ws.Range("E46").Value = "  +1.25%  "

$ws.Range("E47").Value = "  +2.00%  "

$ws.Range("D48").Value = "'0.0242"
$ws.Range("E48").Value = "  +3.88%  "

$ws.Range("D49").Value = "'10.68"
$ws.Range("E49").Value = "  +0.51%  "

$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "'4.82"
$ws.Range("E50").Value = "  +6.90%  "

$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "'19.44"
$ws.Range("E51").Value = "  +4.83%  "
